$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F1").Value = "sd_PAR"

$values = @(
    0.474598139063171,
    396.284463627498,
    347.39216649403,
    223.179808834595,
    0.0372677996249965,
    345.606601886745,
    335.451199212771,
    57.0014081034743,
    0.0372677996249965,
    322.272102408929,
    246.535919304619,
    152.849375366198,
    3.635242764271,
    285.629184228828,
    271.060039686078,
    149.695680087256
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
